$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D:E columns as Text first so numeric-looking values (prices, percentages)
# are stored as text strings, matching the source data which are all inline strings.
$ws.Range("D2:E48").NumberFormat = "@"

$ws.Range("D2").Value = "244.29"
$ws.Range("E2").Value = "-0.72%"

$ws.Range("D3").Value = "27.17"
$ws.Range("E3").Value = "3.89%"

$ws.Range("D4").Value = "5.154"
$ws.Range("E4").Value = "1.08%"

$ws.Range("D5").Value = "0.05632"
$ws.Range("E5").Value = "0.37%"

$ws.Range("D6").Value = "6.474"
$ws.Range("E6").Value = "-0.11%"

$ws.Range("D7").Value = "0.8166"
$ws.Range("E7").Value = "0.61%"

$ws.Range("D8").Value = "0.8329"
$ws.Range("E8").Value = "-1.71%"

$ws.Range("D9").Value = "0.1330"
$ws.Range("E9").Value = "-1.24%"

$ws.Range("D10").Value = "0.06914"
$ws.Range("E10").Value = "-0.77%"

$ws.Range("D11").Value = "0.02935"
$ws.Range("E11").Value = "6.13%"

$ws.Range("D12").Value = "0.09396"
$ws.Range("E12").Value = "-0.09%"

$ws.Range("D13").Value = "0.001508"
$ws.Range("E13").Value = "-0.13%"

$ws.Range("D14").Value = "0.04235"
$ws.Range("E14").Value = "-9.72%"

$ws.Range("D15").Value = "0.0005957"
$ws.Range("E15").Value = "-0.68%"

$ws.Range("D16").Value = "0.006154"
$ws.Range("E16").Value = "0.73%"

$ws.Range("D17").Value = "3.558"
$ws.Range("E17").Value = "0.04%"

$ws.Range("D18").Value = "3.004"
$ws.Range("E18").Value = "-0.57%"

$ws.Range("D19").Value = "2.227"
$ws.Range("E19").Value = "5.12%"

$ws.Range("E20").Value = "-2.20%"

$ws.Range("E21").Value = "-3.31%"

$ws.Range("E22").Value = "-2.12%"

$ws.Range("D23").Value = "3.746"
$ws.Range("E23").Value = "0.08%"

$ws.Range("E24").Value = "-0.08%"

$ws.Range("D25").Value = "0.001224"
$ws.Range("E25").Value = "-1.83%"

$ws.Range("D26").Value = "0.004482"
$ws.Range("E26").Value = "-2.96%"

$ws.Range("E27").Value = "2.04%"

$ws.Range("E28").Value = "-0.45%"

$ws.Range("D40").Value = "0.03652"
$ws.Range("E40").Value = "-0.10%"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1379"
$ws.Range("E41").Value = "1.91%"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002609"
$ws.Range("E42").Value = "-1.92%"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003436"
$ws.Range("E43").Value = "-43.75%"

$ws.Range("D44").Value = "0.008212"
$ws.Range("E44").Value = "-4.85%"

$ws.Range("D45").Value = "0.00005393"
$ws.Range("E45").Value = "1.94%"

$ws.Range("D47").Value = "0.1090"
$ws.Range("E47").Value = "-18.05%"

$ws.Range("D48").Value = "0.002645"
$ws.Range("E48").Value = "29.08%"
